$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: add a new "Prototipação: Figma." bullet before "Outros recursos:
# RoughAnimator." in the Tecnologias e Recursos list.
# ---------------------------------------------------------------------------

$target = $d.Content
$target.Find.Execute("Outros recursos: ") | Out-Null
$targetPara = $target.Paragraphs(1)
$targetIndex = $targetPara.Index

$targetPara.Range.InsertParagraphBefore() | Out-Null

$insertedPara = $d.Paragraphs.Item($targetIndex)
$insertRange = $insertedPara.Range
$insertRange.Collapse(1)

$xml1 = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="4"/>
    </w:numPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Candara" w:hAnsi="Candara"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Candara" w:hAnsi="Candara"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Prototipação: </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Candara" w:hAnsi="Candara"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Figma</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Candara" w:hAnsi="Candara"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
</w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertRange.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Change 2: merge the two runs describing the missing "Envio de
# sugestão/mensagem" feature into a single run, moving the
# <w:lastRenderedPageBreak/> marker to the very start of that run.
# ---------------------------------------------------------------------------

$limFind = $d.Content
$limFind.Find.Execute("O projeto não possui função") | Out-Null
$limPara = $limFind.Paragraphs(1)

$mergedRange = $limPara.Range

$xml2 = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Candara" w:hAnsi="Candara"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r w:rsidRPr="00CF2DA8">
    <w:rPr>
      <w:rFonts w:ascii="Candara" w:hAnsi="Candara"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t xml:space="preserve">O projeto não possui função que recebe dados do usuário, como sugestões e mensagem. Tal limitação deve ser corrigida através da implementação da função &#8220;Envio de sugestão/mensagem&#8221;. Essa funcionalidade deve ser feita através de </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r w:rsidRPr="00CF2DA8">
    <w:rPr>
      <w:rFonts w:ascii="Candara" w:hAnsi="Candara"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>BackEnd</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r w:rsidRPr="00CF2DA8">
    <w:rPr>
      <w:rFonts w:ascii="Candara" w:hAnsi="Candara"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> com Java e </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r w:rsidRPr="00CF2DA8">
    <w:rPr>
      <w:rFonts w:ascii="Candara" w:hAnsi="Candara"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>SpringBoot</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r w:rsidRPr="00CF2DA8">
    <w:rPr>
      <w:rFonts w:ascii="Candara" w:hAnsi="Candara"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
</w:p>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$mergedRange.InsertXML($xml2)
